# Update the "settings" sheet to add two new pseudo-setting rows for
# section1 / section2 (so the generated form emits a _section prompt and
# shows section titles in the table of contents), pushing the existing
# "default" / "hindi" language rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Insert two new blank rows above the current row 5 ("default" language row).
$ws.Range("A5:A6").EntireRow.Insert()

# Row 5: section1 -> display.title "Section 1" / display.title.hindi "धारा 1"
$ws.Range("A5").Value = "section1"
# Row 6: section2 -> display.title "Section 2" / display.title.hindi "धारा 2"
$ws.Range("A6").Value = "section2"
$ws.Range("C6").Value = "Section 2"
$ws.Range("C5").Value = "Section 1"
$ws.Range("D5").Value = "धारा 1"
$ws.Range("D6").Value = "धारा 2"

$ws.Range("C10").Select()
